# Update metadata values on the "Metadata" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 2.0.2
$ws.Range("B3").Value = "2.0.2"

# Date: 2024-04-25T17:24:48+00:00 -> 2025-02-05T10:42:38+00:00
$ws.Range("B8").Value = "2025-02-05T10:42:38+00:00"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# pushing every following row down by one.
$ws.Rows.Item(11).Insert()

# The freshly inserted row doesn't inherit the table's normal cell style,
# so copy formatting down from the row that follows it (now row 12, the
# old "Description" row) before filling in the new values.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
